# Handback status report generation: add a second handed-back file
# (76427fc3-37b4-4a8f-b719-71cf6a2f37cc.md) alongside the existing
# 24659008-... entry, which itself is renamed to 0cb511f5-... and
# given refreshed timestamps.

$wb = $excel.ActiveWorkbook

$oldUuid = "24659008-6d08-43ac-8072-5cfee2563ff9"
$newUuid1 = "0cb511f5-8326-4fb2-92d2-badb413f3755"
$newUuid2 = "76427fc3-37b4-4a8f-b719-71cf6a2f37cc"

$oldZhHash = "013fa165bbd107f6d1205b71446064b63bcc4385"
$newHash1 = "572194832580756001e2705d27cda3b1727b6c80"
$newHash2 = "3723eb50359a2416414c08b22c99cb1ae9df6613"

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newUuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newUuid1.md"
$wsOverview.Range("G2").Value = "2016-08-20 01:04:09"

$wsOverview.Range("A3").Value = "$newUuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newUuid2.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-08-20 01:04:09"

Write-Host "overview done"

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 2 : refresh file name / hash / timestamps for the renamed file
$wsZhCn.Range("A2").Value = "$newUuid1.md"
$wsZhCn.Range("G2").Value = "$newUuid1.$newHash1.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-20 01:04:00"
$wsZhCn.Range("I2").Value = "$newUuid1.md"
$wsZhCn.Range("J2").Value = "$newUuid1.$newHash1.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-20 01:04:27"

# Row 3 : second handed-back file (content duplicate of row 2)
$wsZhCn.Range("A3").Value = "$newUuid2.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "True"
$wsZhCn.Range("G3").Value = "$newUuid2.$newHash2.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-20 01:04:00"
$wsZhCn.Range("I3").Value = "$newUuid2.md"
$wsZhCn.Range("J3").Value = "$newUuid2.$newHash2.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-20 01:04:27"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

Write-Host "zh-cn done"

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 2 : refresh file name / hash / timestamps for the renamed file
$wsDeDe.Range("A2").Value = "$newUuid1.md"
$wsDeDe.Range("G2").Value = "$newUuid1.$newHash1.de-de.xlf"
$wsDeDe.Range("I2").Value = "$newUuid1.md"
$wsDeDe.Range("J2").Value = "$newUuid1.$newHash1.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-20 01:04:33"

# Row 3 : second handed-back file (content duplicate of row 2)
$wsDeDe.Range("A3").Value = "$newUuid2.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "True"
$wsDeDe.Range("G3").Value = "$newUuid2.$newHash2.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-20 01:04:09"
$wsDeDe.Range("I3").Value = "$newUuid2.md"
$wsDeDe.Range("J3").Value = "$newUuid2.$newHash2.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-20 01:04:33"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

Write-Host "de-de done"
